$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the two trailing forecast-origin rows and the trailing BA date column that the
# bugfix removed (the naive forecaster window is now one iteration shorter).
$ws.Range("A23:A24").EntireRow.Delete() | Out-Null
$ws.Range("BA1:BA22").EntireColumn.Delete() | Out-Null

# Header row of forecast-origin date serials (cols B..AZ), shifted one column left
# after bugfix (oldest origin column dropped).
$row1Vals = @(39583,39765,39948,40130,40310,40494,40676,40862,41044,41228,41409,41592,41774,41957,42137,42321,42503,42689,42867,43053,43145,43235,43326,43418,43510,43600,43691,43783,43875,43966,44068,44159,44251,44341,44432,44525,44617,44706,44798,44890,44981,45071,45163,45254,45345,45436,45534,45618,45713,45800,45891)
for ($i = 0; $i -lt $row1Vals.Length; $i++) {
    $ws.Cells.Item(1, $i + 2).Value = $row1Vals[$i]
}

# Row labels in column A (component vintage dates), rows 2..22
$rowLabels = @{
    2 = 39813;
    3 = 40178;
    4 = 40543;
    5 = 40908;
    6 = 41274;
    7 = 41639;
    8 = 42004;
    9 = 42369;
    10 = 42735;
    11 = 43100;
    12 = 43465;
    13 = 43830;
    14 = 44196;
    15 = 44561;
    16 = 44926;
    17 = 45291;
    18 = 45657;
    19 = 46022;
    20 = 46387;
    21 = 46752;
    22 = 47118;
}
foreach ($r in $rowLabels.Keys) {
    $ws.Cells.Item($r, 1).Value = $rowLabels[$r]
}

# Recomputed YoY component forecasts, rows 3..22. Clear the full B:AZ span first (layout
# of the staircase shifts per row) then write the corrected contiguous run.
# row 3
$ws.Range($ws.Cells.Item(3,2), $ws.Cells.Item(3,52)).ClearContents() | Out-Null
$rowVals = @(-1.317619545389281,-1.317619545389281,-1.317619545389281,-1.317619545389281,-1.317619545389281,-1.317619545389281,-1.31761954538927,-1.31761954538927,-1.31761954538927,-1.31761954538927,-1.31761954538927,-1.31761954538927,-1.31761954538927,-1.31761954538927,-1.31761954538927,-1.31761954538927,-1.31761954538927,-1.31761954538927,-1.31761954538927,-1.31761954538927,-1.31761954538927,-1.31761954538927,-1.31761954538927,-1.31761954538927,-1.31761954538927,-1.31761954538927,-1.31761954538927,-1.31761954538927,-1.31761954538927,-1.31761954538927,-1.31761954538927,-1.31761954538927,-1.31761954538927,-1.31761954538927,-1.31761954538927,-1.31761954538927,-1.31761954538927,-1.31761954538927,-1.31761954538927,-1.31761954538927,-1.31761954538927,-1.31761954538927,-1.31761954538927,-1.31761954538927,-1.31761954538927,-1.31761954538927,-1.31761954538927,-1.31761954538927)
for ($i = 0; $i -lt $rowVals.Length; $i++) {
    $ws.Cells.Item(3, 5 + $i).Value = $rowVals[$i]
}

# row 4
$ws.Range($ws.Cells.Item(4,2), $ws.Cells.Item(4,52)).ClearContents() | Out-Null
$rowVals = @(2.771597318554297,2.771597318554297,2.771597318554297,2.771597318554297,2.771597318554297,2.771597318554297,2.771597318554297,2.771597318554297,2.771597318554297,2.771597318554297,2.771597318554297,2.771597318554297,2.771597318554297,2.771597318554297,2.771597318554297,2.771597318554297,2.771597318554297,2.771597318554297,2.771597318554297,2.771597318554297,2.771597318554297,2.771597318554297,2.771597318554297,2.771597318554297,2.771597318554297,2.771597318554297,2.771597318554297,2.771597318554297,2.771597318554297,2.771597318554297,2.771597318554297,2.771597318554297,2.771597318554297,2.771597318554297,2.771597318554297,2.771597318554297,2.771597318554297,2.771597318554297,2.771597318554297,2.771597318554297,2.771597318554297,2.771597318554297,2.771597318554297,2.771597318554297,2.771597318554297,2.771597318554297)
for ($i = 0; $i -lt $rowVals.Length; $i++) {
    $ws.Cells.Item(4, 7 + $i).Value = $rowVals[$i]
}

# row 5
$ws.Range($ws.Cells.Item(5,2), $ws.Cells.Item(5,52)).ClearContents() | Out-Null
$rowVals = @(1.799362536952542,1.799362536952542,1.799362536952542,1.799362536952542,1.799362536952542,1.799362536952542,1.799362536952542,1.799362536952542,1.799362536952542,1.799362536952542,1.799362536952542,1.799362536952542,1.799362536952542,1.799362536952542,1.799362536952542,1.799362536952542,1.799362536952542,1.799362536952542,1.799362536952542,1.799362536952542,1.799362536952542,1.799362536952542,1.799362536952542,1.799362536952542,1.799362536952542,1.799362536952542,1.799362536952542,1.799362536952542,1.799362536952542,1.799362536952542,1.799362536952542,1.799362536952542,1.799362536952542,1.799362536952542,1.799362536952542,1.799362536952542,1.799362536952542,1.799362536952542,1.799362536952542,1.799362536952542,1.799362536952542,1.799362536952542,1.799362536952542,1.799362536952542)
for ($i = 0; $i -lt $rowVals.Length; $i++) {
    $ws.Cells.Item(5, 9 + $i).Value = $rowVals[$i]
}

# row 6
$ws.Range($ws.Cells.Item(6,2), $ws.Cells.Item(6,52)).ClearContents() | Out-Null
$rowVals = @(2.123182427147152,2.123182427147152,2.123182427147152,2.123182427147152,2.123182427147152,2.123182427147152,2.123182427147152,2.123182427147152,2.123182427147152,2.123182427147152,2.123182427147152,2.123182427147152,2.123182427147152,2.123182427147152,2.123182427147152,2.123182427147152,2.123182427147152,2.123182427147152,2.123182427147152,2.123182427147152,2.123182427147152,2.123182427147152,2.123182427147152,2.123182427147152,2.123182427147152,2.123182427147152,2.123182427147152,2.123182427147152,2.123182427147152,2.123182427147152,2.123182427147152,2.123182427147152,2.123182427147152,2.123182427147152,2.123182427147152,2.123182427147152,2.123182427147152,2.123182427147152,2.123182427147152,2.123182427147152,2.123182427147152,2.123182427147152)
for ($i = 0; $i -lt $rowVals.Length; $i++) {
    $ws.Cells.Item(6, 11 + $i).Value = $rowVals[$i]
}

# row 7
$ws.Range($ws.Cells.Item(7,2), $ws.Cells.Item(7,52)).ClearContents() | Out-Null
$rowVals = @(4.792854588620821,4.880442637054072,7.317297369612819,7.317297369612819,7.317297369612819,7.317297369612819,7.317297369612819,7.317297369612819,7.317297369612819,7.317297369612819,7.317297369612819,7.317297369612819,7.317297369612819,7.317297369612819,7.317297369612819,7.317297369612819,7.317297369612819,7.317297369612819,7.317297369612819,7.317297369612819,7.317297369612819,7.317297369612819,7.317297369612819,7.317297369612819,7.317297369612819,7.317297369612819,7.317297369612819,7.317297369612819,7.317297369612819,7.317297369612819,7.317297369612819,7.317297369612819,7.317297369612819,7.317297369612819,7.317297369612819,7.317297369612819,7.317297369612819,7.317297369612819,7.317297369612819,7.317297369612819,7.317297369612819,7.317297369612819)
for ($i = 0; $i -lt $rowVals.Length; $i++) {
    $ws.Cells.Item(7, 11 + $i).Value = $rowVals[$i]
}

# row 8
$ws.Range($ws.Cells.Item(8,2), $ws.Cells.Item(8,52)).ClearContents() | Out-Null
$rowVals = @(1.616393216762324,1.641301872652501,7.239454936865775,5.941867202078877,4.260319658857736,4.260319658857736,4.260319658857736,4.260319658857736,4.260319658857736,4.260319658857736,4.260319658857736,4.260319658857736,4.260319658857736,4.260319658857736,4.260319658857736,4.260319658857736,4.260319658857736,4.260319658857736,4.260319658857736,4.260319658857736,4.260319658857736,4.260319658857736,4.260319658857736,4.260319658857736,4.260319658857736,4.260319658857736,4.260319658857736,4.260319658857736,4.260319658857736,4.260319658857736,4.260319658857736,4.260319658857736,4.260319658857736,4.260319658857736,4.260319658857736,4.260319658857736,4.260319658857736,4.260319658857736,4.260319658857736,4.260319658857736,4.260319658857736,4.260319658857736)
for ($i = 0; $i -lt $rowVals.Length; $i++) {
    $ws.Cells.Item(8, 11 + $i).Value = $rowVals[$i]
}

# row 9
$ws.Range($ws.Cells.Item(9,2), $ws.Cells.Item(9,52)).ClearContents() | Out-Null
$rowVals = @(1.577140242525665,3.107753298997817,2.672847571394987,0.4361429468412448,0.292749233164491,-0.05262415810141086,-0.05262415810141086,-0.05262415810141086,-0.05262415810141086,-0.05262415810141086,-0.05262415810141086,-0.05262415810141086,-0.05262415810141086,-0.05262415810141086,-0.05262415810141086,-0.05262415810141086,-0.05262415810141086,-0.05262415810141086,-0.05262415810141086,-0.05262415810141086,-0.05262415810141086,-0.05262415810141086,-0.05262415810141086,-0.05262415810141086,-0.05262415810141086,-0.05262415810141086,-0.05262415810141086,-0.05262415810141086,-0.05262415810141086,-0.05262415810141086,-0.05262415810141086,-0.05262415810141086,-0.05262415810141086,-0.05262415810141086,-0.05262415810141086,-0.05262415810141086,-0.05262415810141086,-0.05262415810141086,-0.05262415810141086,-0.05262415810141086,-0.05262415810141086)
for ($i = 0; $i -lt $rowVals.Length; $i++) {
    $ws.Cells.Item(9, 12 + $i).Value = $rowVals[$i]
}

# row 10
$ws.Range($ws.Cells.Item(10,2), $ws.Cells.Item(10,52)).ClearContents() | Out-Null
$rowVals = @(2.668903200194506,2.107195049295729,2.181874035977249,1.392195163617171,2.032207428223742,1.459778471779982,1.459778471779982,1.459778471779982,1.459778471779982,1.459778471779982,1.459778471779982,1.459778471779982,1.459778471779982,1.459778471779982,1.459778471779982,1.459778471779982,1.459778471779982,1.459778471779982,1.459778471779982,1.459778471779982,1.459778471779982,1.459778471779982,1.459778471779982,1.459778471779982,1.459778471779982,1.459778471779982,1.459778471779982,1.459778471779982,1.459778471779982,1.459778471779982,1.459778471779982,1.459778471779982,1.459778471779982,1.459778471779982,1.459778471779982,1.459778471779982,1.459778471779982,1.459778471779982,1.459778471779982)
for ($i = 0; $i -lt $rowVals.Length; $i++) {
    $ws.Cells.Item(10, 14 + $i).Value = $rowVals[$i]
}

# row 11
$ws.Range($ws.Cells.Item(11,2), $ws.Cells.Item(11,52)).ClearContents() | Out-Null
$rowVals = @(2.116805664817756,1.952530827880228,2.13692496326825,1.493220091771108,2.775332754349846,3.002208343813528,3.002208343813528,3.002208343813528,3.002208343813528,3.002208343813528,3.002208343813528,3.002208343813528,3.002208343813528,3.002208343813528,3.002208343813528,3.002208343813528,3.002208343813528,3.002208343813528,3.002208343813528,3.002208343813528,3.002208343813528,3.002208343813528,3.002208343813528,3.002208343813528,3.002208343813528,3.002208343813528,3.002208343813528,3.002208343813528,3.002208343813528,3.002208343813528,3.002208343813528,3.002208343813528,3.002208343813528,3.002208343813528,3.002208343813528,3.002208343813528,3.002208343813528)
for ($i = 0; $i -lt $rowVals.Length; $i++) {
    $ws.Cells.Item(11, 16 + $i).Value = $rowVals[$i]
}

# row 12
$ws.Range($ws.Cells.Item(12,2), $ws.Cells.Item(12,52)).ClearContents() | Out-Null
$rowVals = @(2.064081081388358,1.928696216909276,2.200426660963761,2.622364272988187,3.2651197821016,3.565025829754953,4.01493878081518,4.020433260014977,4.020433260014977,4.020433260014977,4.020433260014977,4.020433260014977,4.020433260014977,4.020433260014977,4.020433260014977,4.020433260014977,4.020433260014977,4.020433260014977,4.020433260014977,4.020433260014977,4.020433260014977,4.020433260014977,4.020433260014977,4.020433260014977,4.020433260014977,4.020433260014977,4.020433260014977,4.020433260014977,4.020433260014977,4.020433260014977,4.020433260014977,4.020433260014977,4.020433260014977,4.020433260014977,4.020433260014977)
for ($i = 0; $i -lt $rowVals.Length; $i++) {
    $ws.Cells.Item(12, 18 + $i).Value = $rowVals[$i]
}

# row 13
$ws.Range($ws.Cells.Item(13,2), $ws.Cells.Item(13,52)).ClearContents() | Out-Null
$rowVals = @(2.083217522782399,2.156473785802171,2.271936475508851,2.446228176258058,3.091110147865495,3.113086948791377,3.283136334808323,3.444206290325491,3.479628752085517,3.53224976671227,3.53224976671227,3.53224976671227,3.53224976671227,3.53224976671227,3.53224976671227,3.53224976671227,3.53224976671227,3.53224976671227,3.53224976671227,3.53224976671227,3.53224976671227,3.53224976671227,3.53224976671227,3.53224976671227,3.53224976671227,3.53224976671227,3.53224976671227,3.53224976671227,3.53224976671227,3.53224976671227,3.53224976671227,3.53224976671227,3.53224976671227)
for ($i = 0; $i -lt $rowVals.Length; $i++) {
    $ws.Cells.Item(13, 20 + $i).Value = $rowVals[$i]
}

# row 14
$ws.Range($ws.Cells.Item(14,2), $ws.Cells.Item(14,52)).ClearContents() | Out-Null
$rowVals = @(2.260904903527239,2.369764785923656,2.373113736336396,2.403408536719187,2.496958452261078,2.54748094003614,2.757421718286168,2.31260691849986,2.667234932970275,-0.985458715495402,-0.985458715495402,-0.985458715495402,-0.985458715495402,-0.985458715495402,-0.985458715495402,-0.985458715495402,-0.985458715495402,-0.985458715495402,-0.985458715495402,-0.985458715495402,-0.985458715495402,-0.985458715495402,-0.985458715495402,-0.985458715495402,-0.985458715495402,-0.985458715495402,-0.985458715495402,-0.985458715495402,-0.985458715495402,-0.985458715495402)
for ($i = 0; $i -lt $rowVals.Length; $i++) {
    $ws.Cells.Item(14, 23 + $i).Value = $rowVals[$i]
}

# row 15
$ws.Range($ws.Cells.Item(15,2), $ws.Cells.Item(15,52)).ClearContents() | Out-Null
$rowVals = @(2.400584622094115,2.408484542873435,2.438756968913824,2.363182008239928,2.567662999186382,-1.999977666418695,-1.373617952268746,-4.853362183897836,-4.511102905979703,-4.365687260408224,-4.365687260408224,-4.365687260408224,-4.365687260408224,-4.365687260408224,-4.365687260408224,-4.365687260408224,-4.365687260408224,-4.365687260408224,-4.365687260408224,-4.365687260408224,-4.365687260408224,-4.365687260408224,-4.365687260408224,-4.365687260408224,-4.365687260408224,-4.365687260408224)
for ($i = 0; $i -lt $rowVals.Length; $i++) {
    $ws.Cells.Item(15, 27 + $i).Value = $rowVals[$i]
}

# row 16
$ws.Range($ws.Cells.Item(16,2), $ws.Cells.Item(16,52)).ClearContents() | Out-Null
$rowVals = @(2.433771721344358,2.229247642542287,2.358999556713859,0.6599234717970859,0.9582724917052587,1.582150300418306,1.324283050325015,2.700663803921799,1.386772772629241,1.27347919322387,1.27347919322387,1.27347919322387,1.27347919322387,1.27347919322387,1.27347919322387,1.27347919322387,1.27347919322387,1.27347919322387,1.27347919322387,1.27347919322387,1.27347919322387,1.27347919322387)
for ($i = 0; $i -lt $rowVals.Length; $i++) {
    $ws.Cells.Item(16, 31 + $i).Value = $rowVals[$i]
}

# row 17
$ws.Range($ws.Cells.Item(17,2), $ws.Cells.Item(17,52)).ClearContents() | Out-Null
$rowVals = @(2.030062065969385,1.907023604270441,2.050710065226147,2.080813569776696,2.813733698830267,0.8813242377093244,0.4287685802702779,-0.5025420863900898,-0.185315122156382,-0.9537175292835154,-1.339436245206127,-1.339436245206127,-1.339436245206127,-1.339436245206127,-1.339436245206127,-1.339436245206127,-1.339436245206127,-1.339436245206127,-1.339436245206127)
for ($i = 0; $i -lt $rowVals.Length; $i++) {
    $ws.Cells.Item(17, 34 + $i).Value = $rowVals[$i]
}

# row 18
$ws.Range($ws.Cells.Item(18,2), $ws.Cells.Item(18,52)).ClearContents() | Out-Null
$rowVals = @(2.331302864153018,1.729606009313156,1.632639074358488,1.517427911430702,1.715279842342743,0.9049225073274991,-0.6598239038080322,-2.371854438773213,-2.5174493871855,-3.303819519576723,-3.451527003230626,-3.451527003230626,-3.451527003230626,-3.451527003230626,-3.451527003230626)
for ($i = 0; $i -lt $rowVals.Length; $i++) {
    $ws.Cells.Item(18, 38 + $i).Value = $rowVals[$i]
}

# row 19
$ws.Range($ws.Cells.Item(19,2), $ws.Cells.Item(19,52)).ClearContents() | Out-Null
$rowVals = @(1.65517688847745,1.853878336645232,1.635646651659406,1.196659994305449,1.089612584330668,-0.3655818470008065,-1.070961900287937,-1.436963918858969,-1.652703512303566,-1.704805397136089,-1.75044229618867)
for ($i = 0; $i -lt $rowVals.Length; $i++) {
    $ws.Cells.Item(19, 42 + $i).Value = $rowVals[$i]
}

# row 20
$ws.Range($ws.Cells.Item(20,2), $ws.Cells.Item(20,52)).ClearContents() | Out-Null
$rowVals = @(1.412018118185099,0.9129647504396488,0.6686113921051895,0.5903916006237075,0.4469810487905734,0.1460701281005727,-0.1564272439962933)
for ($i = 0; $i -lt $rowVals.Length; $i++) {
    $ws.Cells.Item(20, 46 + $i).Value = $rowVals[$i]
}

# row 21
$ws.Range($ws.Cells.Item(21,2), $ws.Cells.Item(21,52)).ClearContents() | Out-Null
$rowVals = @(0.6478620066787233,0.2977127898601273,0.03402167163586434)
for ($i = 0; $i -lt $rowVals.Length; $i++) {
    $ws.Cells.Item(21, 50 + $i).Value = $rowVals[$i]
}

# row 22
$ws.Range($ws.Cells.Item(22,2), $ws.Cells.Item(22,52)).ClearContents() | Out-Null

Write-Output "edit complete"